$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "queen" -> "pencil" (reuse of row 6), plus a brand new word "powerbank" in row 7
$ws.Range("A6").Value = "pencil"
$ws.Range("A7").Value = "powerbank"

# Mirror the author's final on-screen selection (rows 6:9, anchored near the new entries)
$ws.Range("A6:A9").Select() | Out-Null
